$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 21:20"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 180789
$ws.Range("C4").Value = 17001
$ws.Range("E4").Value = 170968
$ws.Range("G4").Value = 439
$ws.Range("H4").Value = 3580

# Row 8 - Alemania
$ws.Range("B8").Value = 70985
$ws.Range("C8").Value = 4100
$ws.Range("E8").Value = 54479

# Row 20 - Israel
$ws.Range("B20").Value = 5358
$ws.Range("C20").Value = 663
$ws.Range("D20").Value = 224
$ws.Range("E20").Value = 5114
$ws.Range("F20").Value = 117

# Row 21 - Brasil
$ws.Range("B21").Value = 4725
$ws.Range("C21").Value = 95
$ws.Range("E21").Value = 4430

# Row 22 - Noruega
$ws.Range("B22").Value = 4630
$ws.Range("C22").Value = 185
$ws.Range("E22").Value = 4578

# Row 26 - Irlanda
$ws.Range("B26").Value = 3235
$ws.Range("C26").Value = 325
$ws.Range("E26").Value = 3159
$ws.Range("G26").Value = 17
$ws.Range("H26").Value = 71

# Row 32 - Rumania
$ws.Range("E32").Value = 1943
$ws.Range("G32").Value = 17
$ws.Range("H32").Value = 82

# Row 72 - Bosnia y Herzegovina
$ws.Range("B72").Value = 418
$ws.Range("C72").Value = 50
$ws.Range("E72").Value = 389

# Row 135 - Polinesia Francesa
$ws.Range("F135").Value = 1
